$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 207 (a new weekly price observation), pushing the
# existing rows 207-222 down to 208-223. Duplicate row 207's contents first
# so the static columns (Mercado, Region, Codreg, Categoria, etc.) carry
# over, then overwrite the columns that differ for the new observation
# (Fecha, Volumen, Precio minimo/maximo/promedio, Precio $/Kg).
$ws.Rows.Item(207).Copy()
$ws.Rows.Item(207).Insert()

$ws.Cells.Item(207, 4).Value = 44610
$ws.Cells.Item(207, 10).Value = 220
$ws.Cells.Item(207, 11).Value = 21000
$ws.Cells.Item(207, 12).Value = 21000
$ws.Cells.Item(207, 13).Value = 21000
$ws.Cells.Item(207, 16).Value = 2100
